# Refresh cryptos list price/volume snapshot (GitHub Actions scheduled update).
# Source diff only touches column D (Price) and column E (Volume(1h)) text values
# for the existing data rows (2-51); everything else is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.869.47"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "3.543.58"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.19"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.81"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "3.541.33"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "4.143.06"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.03"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "3.544.62"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "67.656.40"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.25"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.69"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.74"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.624"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.28"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("E25").Value = "  +5.95%  "
$ws.Range("D26").Value = "3.685.30"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.70"
$ws.Range("E29").Value = "  +4.65%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  +6.60%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.75"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.23"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "3.531.34"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.06"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "176.24"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.887"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.06"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.58"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.994"
$ws.Range("E51").Value = "  -4.18%  "
